$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, pushing the existing rows 80-92 down to 81-93.
$ws.Rows.Item(80).Insert()

# The new blank row 80 inherited the date style (s="2") on D80 from the insert
# (Excel's default "format from row above" behavior). Populate it with the new
# record's data; columns A, B, C, E, F, G, H, I, R carry the same constant
# values used throughout this block of rows.
$ws.Cells.Item(80, 1).Value = 8
$ws.Cells.Item(80, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(80, 3).Value = "Coquimbo"
$ws.Cells.Item(80, 4).Value = 44722
$ws.Cells.Item(80, 5).Value = 4
$ws.Cells.Item(80, 6).Value = 100112052
$ws.Cells.Item(80, 7).Value = "Albahaca"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 1100
$ws.Cells.Item(80, 11).Value = 3500
$ws.Cells.Item(80, 12).Value = 4000
$ws.Cells.Item(80, 13).Value = 3750
$ws.Cells.Item(80, 14).Value = "`$/paquete"
$ws.Cells.Item(80, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(80, 16).Value = 3750
$ws.Cells.Item(80, 17).Value = 1
$ws.Cells.Item(80, 18).Value = "Hortaliza"
